$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free direct assignment. For D-column values that parse as plain
# numbers, force a Text number format first (then restore to Normal style)
# so Excel keeps the original zero-padded / fixed-decimal string instead of
# silently converting it to a numeric value.

$ws.Range('D2').Value = '68.887.27'
$ws.Range('E2').Value = '  -3.72%  '
$ws.Range('D3').Value = '3.489.93'
$ws.Range('E3').Value = '  -5.85%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.32'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.604'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.01%  '
$ws.Range('D8').Value = '3.483.18'
$ws.Range('E8').Value = '  -5.63%  '
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.187'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.47'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.578'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '46.10'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000271'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.65%  '
$ws.Range('D15').Value = '4.051.33'
$ws.Range('E15').Value = '  -5.90%  '
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.47'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -6.09%  '
$ws.Range('B17').Value = 'BitcoinCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '622.21'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -9.00%  '
$ws.Range('D18').Value = '68.832.15'
$ws.Range('E18').Value = '  -3.99%  '
$ws.Range('D19').Value = '3.478.56'
$ws.Range('E19').Value = '  -6.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.121'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.877'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '15.74'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '96.57'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.77'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.83%  '
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('E28').Value = '  -8.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.22'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -11.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.39'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -8.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.14'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -8.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.43'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -8.25%  '
$ws.Range('E33').Value = '  -9.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.91'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '623.46'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.68'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.76%  '
$ws.Range('E37').Value = '  -6.11%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.41'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -16.89%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '56.34'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0443'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.57%  '
$ws.Range('E42').Value = '  -6.90%  '
$ws.Range('D43').Value = '3.331.76'
$ws.Range('E43').Value = '  -9.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.324'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '32.45'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -8.60%  '
$ws.Range('D46').Value = '0.0₃0683'
$ws.Range('E46').Value = '  -11.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.54'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -9.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.76'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.129'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '131.04'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.50%  '
$ws.Range('E51').Value = '  +13.79%  '
